# The MTG table for rows 40-65 had its "symbol" value shifted one/two
# columns too far right (D/E instead of B/C), with a stray, always-empty
# column (F) separating it from the numeric "length/diameter" data in G.
# Realign it to match the layout used by the reference table above
# (rows 13-16): symbol in B/C, then a single spacer column (D), then the
# numeric data starting at G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$moves = @(
  @("D40","B40"), @("F40","D40"),
  @("D41","B41"), @("F41","D41"),
  @("E42","C42"), @("F42","D42"),
  @("E43","C43"), @("F43","D43"),
  @("D44","B44"), @("F44","D44"),
  @("E45","C45"), @("F45","D45"),
  @("E46","C46"), @("F46","D46"),
  @("D47","B47"), @("F47","D47"),
  @("E48","C48"), @("F48","D48"),
  @("E49","C49"), @("F49","D49"),
  @("D50","B50"), @("F50","D50"),
  @("E51","C51"), @("F51","D51"),
  @("E52","C52"), @("F52","D52"),
  @("D53","B53"), @("F53","D53"),
  @("E54","C54"), @("F54","D54"),
  @("E55","C55"), @("F55","D55"),
  @("D56","B56"), @("F56","D56"),
  @("E57","C57"), @("F57","D57"),
  @("E58","C58"), @("F58","D58"),
  @("D59","B59"), @("F59","D59"),
  @("E60","C60"), @("F60","D60"),
  @("E61","C61"), @("F61","D61"),
  @("D62","B62"), @("F62","D62"),
  @("E63","C63"), @("F63","D63"),
  @("E64","C64"), @("F64","D64"),
  @("D65","B65"), @("F65","D65")
)

foreach ($m in $moves) {
  $src = $m[0]
  $dst = $m[1]
  $ws.Range($src).Copy($ws.Range($dst))
  $ws.Range($src).Clear()
}

# Match the saved selection / scroll state on the sheet.
$ws.Range("E38").Select()
